$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = "0e47dd16-3ece-4c4b-ba97-3101347d82c1"
$ws.Range("B6").Value = "In"
$ws.Range("C6").Value = "One Face"
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 10
# Leading apostrophe forces literal text so the date-shaped string isn't
# auto-converted into an Excel date serial number.
$ws.Range("F6").Value = "'2024-09-23"
$ws.Range("G6").Value = "21:21:45"

# Row 7
$ws.Range("A7").Value = "83cd03f2-981e-4881-9a90-788063c9152b"
$ws.Range("B7").Value = "Waste"
$ws.Range("C7").Value = "paper"
$ws.Range("D7").Value = 100
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "'2024-09-23"
$ws.Range("G7").Value = "21:23:05"
